$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add new worksheet "cohort_fake" after the existing sheet ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "cohort_fake"

# --- Populate the new sheet with the cohort/graduated data table ---
$data = @(
    @("Value", "graduated", 30, 36, 24, 29, 18, 23, 12, 17, 6, 11),
    @(5, 1, 1, 1, 0, 0, 0, 0, 0, 0, 0, 0),
    @(4, 1, 0, 0, 1, 1, 0, 0, 0, 0, 0, 0),
    @(3, 1, 0, 0, 0, 0, 1, 1, 0, 0, 0, 0),
    @(2, 1, 0, 0, 0, 0, 0, 0, 1, 1, 0, 0),
    @(1, 1, 0, 0, 0, 0, 0, 0, 0, 0, 1, 1),
    @(1, 1, 0, 0, 0, 0, 0, 0, 0, 0, 1, 0)
)

for ($r = 1; $r -le $data.Count; $r++) {
    $row = $data[$r - 1]
    for ($c = 1; $c -le $row.Count; $c++) {
        $ws2.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}

# Approximate the column width auto-computed by Excel for column C on this
# sheet (closest value representable through this runtime's width rounding).
$ws2.Columns.Item(3).ColumnWidth = 11.8

$ws2.Range("G11").Select()

# --- Update the original sheet: add a new row with a blank-space marker cell ---
$ws1.Range("F16").Value = " "
$ws1.Range("K20").Select()
